# Natmi following Dr Hou advice
# Update LR-pair stats (Mdk-Itga6) on Sheet1: ligand/receptor expressing-cell
# counts change from 1 to 3, with corresponding recomputed expression /
# specificity / edge-weight values for every data row (rows 2-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 1.324023666666666
$ws.Cells.Item(2, 8).Value = 3.972071
$ws.Cells.Item(2, 9).Value = 0.01518042398701374
$ws.Cells.Item(2, 10).Value = 0.01518042398701374
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 192.8285726666667
$ws.Cells.Item(2, 14).Value = 578.485718
$ws.Cells.Item(2, 15).Value = 0.7801188850698786
$ws.Cells.Item(2, 16).Value = 0.7801188850698786
$ws.Cells.Item(2, 17).Value = 255.3095938202198
$ws.Cells.Item(2, 18).Value = 2297.786344381978
$ws.Cells.Item(2, 19).Value = 0.0118425354356372
$ws.Cells.Item(2, 20).Value = 0.0118425354356372
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 1.324023666666666
$ws.Cells.Item(3, 8).Value = 3.972071
$ws.Cells.Item(3, 9).Value = 0.01518042398701374
$ws.Cells.Item(3, 10).Value = 0.01518042398701374
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.4209206666666667
$ws.Cells.Item(3, 14).Value = 1.262762
$ws.Cells.Item(3, 15).Value = 0.001702901995496819
$ws.Cells.Item(3, 16).Value = 0.001702901995496819
$ws.Cells.Item(3, 17).Value = 0.5573089244557777
$ws.Cells.Item(3, 18).Value = 5.015780320102
$ws.Cells.Item(3, 19).Value = 0.00002585077429997347
$ws.Cells.Item(3, 20).Value = 0.00002585077429997347
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 1.324023666666666
$ws.Cells.Item(4, 8).Value = 3.972071
$ws.Cells.Item(4, 9).Value = 0.01518042398701374
$ws.Cells.Item(4, 10).Value = 0.01518042398701374
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 45.70525533333333
$ws.Cells.Item(4, 14).Value = 137.115766
$ws.Cells.Item(4, 15).Value = 0.184907933193646
$ws.Cells.Item(4, 16).Value = 0.184907933193646
$ws.Cells.Item(4, 17).Value = 60.51483975237621
$ws.Cells.Item(4, 18).Value = 544.633557771386
$ws.Cells.Item(4, 19).Value = 0.002806980824441958
$ws.Cells.Item(4, 20).Value = 0.002806980824441958
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 1.324023666666666
$ws.Cells.Item(5, 8).Value = 3.972071
$ws.Cells.Item(5, 9).Value = 0.01518042398701374
$ws.Cells.Item(5, 10).Value = 0.01518042398701374
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 8.223696
$ws.Cells.Item(5, 14).Value = 24.671088
$ws.Cells.Item(5, 15).Value = 0.0332702797409786
$ws.Cells.Item(5, 16).Value = 0.0332702797409786
$ws.Cells.Item(5, 17).Value = 10.888368131472
$ws.Cells.Item(5, 18).Value = 97.995313183248
$ws.Cells.Item(5, 19).Value = 0.0005050569526346088
$ws.Cells.Item(5, 20).Value = 0.0005050569526346087
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 81.17653533333333
$ws.Cells.Item(6, 8).Value = 243.529606
$ws.Cells.Item(6, 9).Value = 0.9307191821270077
$ws.Cells.Item(6, 10).Value = 0.9307191821270075
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 192.8285726666667
$ws.Cells.Item(6, 14).Value = 578.485718
$ws.Cells.Item(6, 15).Value = 0.7801188850698786
$ws.Cells.Item(6, 16).Value = 0.7801188850698786
$ws.Cells.Item(6, 17).Value = 15653.1554423519
$ws.Cells.Item(6, 18).Value = 140878.3989811671
$ws.Cells.Item(6, 19).Value = 0.7260716106740704
$ws.Cells.Item(6, 20).Value = 0.7260716106740704
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 81.17653533333333
$ws.Cells.Item(7, 8).Value = 243.529606
$ws.Cells.Item(7, 9).Value = 0.9307191821270077
$ws.Cells.Item(7, 10).Value = 0.9307191821270075
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.4209206666666667
$ws.Cells.Item(7, 14).Value = 1.262762
$ws.Cells.Item(7, 15).Value = 0.001702901995496819
$ws.Cells.Item(7, 16).Value = 0.001702901995496819
$ws.Cells.Item(7, 17).Value = 34.16888137019689
$ws.Cells.Item(7, 18).Value = 307.5199323317721
$ws.Cells.Item(7, 19).Value = 0.001584923552491249
$ws.Cells.Item(7, 20).Value = 0.001584923552491248
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 81.17653533333333
$ws.Cells.Item(8, 8).Value = 243.529606
$ws.Cells.Item(8, 9).Value = 0.9307191821270077
$ws.Cells.Item(8, 10).Value = 0.9307191821270075
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 45.70525533333333
$ws.Cells.Item(8, 14).Value = 137.115766
$ws.Cells.Item(8, 15).Value = 0.184907933193646
$ws.Cells.Item(8, 16).Value = 0.184907933193646
$ws.Cells.Item(8, 17).Value = 3710.194274485355
$ws.Cells.Item(8, 18).Value = 33391.7484703682
$ws.Cells.Item(8, 19).Value = 0.1720973603507856
$ws.Cells.Item(8, 20).Value = 0.1720973603507856
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 81.17653533333333
$ws.Cells.Item(9, 8).Value = 243.529606
$ws.Cells.Item(9, 9).Value = 0.9307191821270077
$ws.Cells.Item(9, 10).Value = 0.9307191821270075
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 8.223696
$ws.Cells.Item(9, 14).Value = 24.671088
$ws.Cells.Item(9, 15).Value = 0.0332702797409786
$ws.Cells.Item(9, 16).Value = 0.0332702797409786
$ws.Cells.Item(9, 17).Value = 667.5711489145921
$ws.Cells.Item(9, 18).Value = 6008.140340231328
$ws.Cells.Item(9, 19).Value = 0.03096528754966036
$ws.Cells.Item(9, 20).Value = 0.03096528754966035
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 1.192675
$ws.Cells.Item(10, 8).Value = 3.578025
$ws.Cells.Item(10, 9).Value = 0.0136744626508778
$ws.Cells.Item(10, 10).Value = 0.0136744626508778
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 192.8285726666667
$ws.Cells.Item(10, 14).Value = 578.485718
$ws.Cells.Item(10, 15).Value = 0.7801188850698786
$ws.Cells.Item(10, 16).Value = 0.7801188850698786
$ws.Cells.Item(10, 17).Value = 229.9818179052167
$ws.Cells.Item(10, 18).Value = 2069.83636114695
$ws.Cells.Item(10, 19).Value = 0.01066770655713249
$ws.Cells.Item(10, 20).Value = 0.01066770655713249
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 1.192675
$ws.Cells.Item(11, 8).Value = 3.578025
$ws.Cells.Item(11, 9).Value = 0.0136744626508778
$ws.Cells.Item(11, 10).Value = 0.0136744626508778
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.4209206666666667
$ws.Cells.Item(11, 14).Value = 1.262762
$ws.Cells.Item(11, 15).Value = 0.001702901995496819
$ws.Cells.Item(11, 16).Value = 0.001702901995496819
$ws.Cells.Item(11, 17).Value = 0.5020215561166668
$ws.Cells.Item(11, 18).Value = 4.518194005050001
$ws.Cells.Item(11, 19).Value = 0.00002328626973552653
$ws.Cells.Item(11, 20).Value = 0.00002328626973552653
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 1.192675
$ws.Cells.Item(12, 8).Value = 3.578025
$ws.Cells.Item(12, 9).Value = 0.0136744626508778
$ws.Cells.Item(12, 10).Value = 0.0136744626508778
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 45.70525533333333
$ws.Cells.Item(12, 14).Value = 137.115766
$ws.Cells.Item(12, 15).Value = 0.184907933193646
$ws.Cells.Item(12, 16).Value = 0.184907933193646
$ws.Cells.Item(12, 17).Value = 54.51151540468334
$ws.Cells.Item(12, 18).Value = 490.6036386421501
$ws.Cells.Item(12, 19).Value = 0.002528516626307521
$ws.Cells.Item(12, 20).Value = 0.00252851662630752
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 1.192675
$ws.Cells.Item(13, 8).Value = 3.578025
$ws.Cells.Item(13, 9).Value = 0.0136744626508778
$ws.Cells.Item(13, 10).Value = 0.0136744626508778
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 8.223696
$ws.Cells.Item(13, 14).Value = 24.671088
$ws.Cells.Item(13, 15).Value = 0.0332702797409786
$ws.Cells.Item(13, 16).Value = 0.0332702797409786
$ws.Cells.Item(13, 17).Value = 9.808196626800001
$ws.Cells.Item(13, 18).Value = 88.2737696412
$ws.Cells.Item(13, 19).Value = 0.0004549531977022684
$ws.Cells.Item(13, 20).Value = 0.0004549531977022682
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 3.525915333333334
$ws.Cells.Item(14, 8).Value = 10.577746
$ws.Cells.Item(14, 9).Value = 0.04042593123510095
$ws.Cells.Item(14, 10).Value = 0.04042593123510094
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 192.8285726666667
$ws.Cells.Item(14, 14).Value = 578.485718
$ws.Cells.Item(14, 15).Value = 0.7801188850698786
$ws.Cells.Item(14, 16).Value = 0.7801188850698786
$ws.Cells.Item(14, 17).Value = 679.897221070181
$ws.Cells.Item(14, 18).Value = 6119.074989631629
$ws.Cells.Item(14, 19).Value = 0.03153703240303853
$ws.Cells.Item(14, 20).Value = 0.03153703240303853
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 3.525915333333334
$ws.Cells.Item(15, 8).Value = 10.577746
$ws.Cells.Item(15, 9).Value = 0.04042593123510095
$ws.Cells.Item(15, 10).Value = 0.04042593123510094
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.4209206666666667
$ws.Cells.Item(15, 14).Value = 1.262762
$ws.Cells.Item(15, 15).Value = 0.001702901995496819
$ws.Cells.Item(15, 16).Value = 0.001702901995496819
$ws.Cells.Item(15, 17).Value = 1.484130632716889
$ws.Cells.Item(15, 18).Value = 13.357175694452
$ws.Cells.Item(15, 19).Value = 0.00006884139897007059
$ws.Cells.Item(15, 20).Value = 0.00006884139897007058
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 3.525915333333334
$ws.Cells.Item(16, 8).Value = 10.577746
$ws.Cells.Item(16, 9).Value = 0.04042593123510095
$ws.Cells.Item(16, 10).Value = 0.04042593123510094
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 45.70525533333333
$ws.Cells.Item(16, 14).Value = 137.115766
$ws.Cells.Item(16, 15).Value = 0.184907933193646
$ws.Cells.Item(16, 16).Value = 0.184907933193646
$ws.Cells.Item(16, 17).Value = 161.1528605937151
$ws.Cells.Item(16, 18).Value = 1450.375745343436
$ws.Cells.Item(16, 19).Value = 0.007475075392110975
$ws.Cells.Item(16, 20).Value = 0.007475075392110974
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 3.525915333333334
$ws.Cells.Item(17, 8).Value = 10.577746
$ws.Cells.Item(17, 9).Value = 0.04042593123510095
$ws.Cells.Item(17, 10).Value = 0.04042593123510094
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 8.223696
$ws.Cells.Item(17, 14).Value = 24.671088
$ws.Cells.Item(17, 15).Value = 0.0332702797409786
$ws.Cells.Item(17, 16).Value = 0.0332702797409786
$ws.Cells.Item(17, 17).Value = 28.99605582307201
$ws.Cells.Item(17, 18).Value = 260.964502407648
$ws.Cells.Item(17, 19).Value = 0.001344982040981373
$ws.Cells.Item(17, 20).Value = 0.001344982040981373
